$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Command" rectangle (UML class box) that needs the {abstract} stereotype
# label added above its class-name line, per the "Fix DG UML bug" commit.
$shp = $s.Shapes.Item(4)

# Resize / reposition the box to make room for the extra line of text.
# (Tiny epsilons nudge the point->EMU conversion to the exact target value,
# avoiding float round-trip truncation for values like 330.3 / 39.05.)
$shp.Left = 330.30001
$shp.Top = 92
$shp.Width = 125.15
$shp.Height = 39.05001

# Insert the new "{abstract}" paragraph above the existing "Command" text,
# inheriting the run formatting already on that text (bold italic 14pt, white).
$tr = $shp.TextFrame.TextRange
$tr.InsertBefore("{abstract}" + [char]13) | Out-Null

$abstractPara = $tr.Paragraphs(1,1)
$abstractPara.Font.Bold = $true
$abstractPara.Font.Italic = $true
$abstractPara.Font.Size = 14
